$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename HT goals headers (shared strings) ---
$ws.Range("I1").Value = "HTHG"
$ws.Range("J1").Value = "HTAG"

# --- 2. Swap duplicate/mismatched match rows (id stays sequential; underlying
#      match data for each pair was swapped between the two rows) ---

# --- Swap row 71 and row 72 (columns B, E:AD); A/C/D unchanged ---
$ws.Range("B71").Value = 6139072
$ws.Range("B72").Value = 6139071
$ws.Range("E71").Value = 'JK Tammeka Tartu'
$ws.Range("E72").Value = 'Parnu JK Vaprus'
$ws.Range("F71").Value = 'FC Flora Tallinn'
$ws.Range("F72").Value = 'JK Trans Narva'
$ws.Range("G71").Value = 1
$ws.Range("G72").Value = 3
$ws.Range("H71").Value = 2
$ws.Range("H72").Value = 2
$ws.Range("I71").Value = 1
$ws.Range("I72").Value = 1
$ws.Range("J71").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K71").Value = 'A'
$ws.Range("K72").Value = 'H'
$ws.Range("L71").Value = 9
$ws.Range("L72").Value = 2.4
$ws.Range("M71").Value = 7
$ws.Range("M72").Value = 3.2
$ws.Range("N71").Value = 1.166
$ws.Range("N72").Value = 2.6
$ws.Range("O71").Value = 7
$ws.Range("O72").Value = 3
$ws.Range("P71").Value = 6
$ws.Range("P72").Value = 3.25
$ws.Range("Q71").Value = 1.25
$ws.Range("Q72").Value = 2.2
$ws.Range("R71").Value = 1.75
$ws.Range("R72").Value = 0.25
$ws.Range("S71").Value = 1.9
$ws.Range("S72").Value = 1.825
$ws.Range("T71").Value = 1.9
$ws.Range("T72").Value = 1.975
$ws.Range("U71").Value = 3
$ws.Range("U72").Value = 2.5
$ws.Range("V71").Value = 1.95
$ws.Range("V72").Value = 1.875
$ws.Range("W71").Value = 1.85
$ws.Range("W72").Value = 1.925
$ws.Range("X71").Value = -1
$ws.Range("X72").Value = 2
$ws.Range("Y71").Value = -1
$ws.Range("Y72").Value = -1
$ws.Range("Z71").Value = 0.25
$ws.Range("Z72").Value = -1
$ws.Range("AA71").Value = 0.8999999999999999
$ws.Range("AA72").Value = 0.825
$ws.Range("AB71").Value = -1
$ws.Range("AB72").Value = -1
$ws.Range("AC71").Value = 0
$ws.Range("AC72").Value = 0.875
$ws.Range("AD71").Value = 0
$ws.Range("AD72").Value = -1

# --- Swap row 88 and row 89 (columns B, E:AD); A/C/D unchanged ---
$ws.Range("B88").Value = 6376947
$ws.Range("B89").Value = 6376945
$ws.Range("E88").Value = 'JK Tammeka Tartu'
$ws.Range("E89").Value = 'Parnu JK Vaprus'
$ws.Range("F88").Value = 'JK Tallinna Kalev'
$ws.Range("F89").Value = 'Harju JK Laagri'
$ws.Range("G88").Value = 2
$ws.Range("G89").Value = 0
$ws.Range("H88").Value = 7
$ws.Range("H89").Value = 0
$ws.Range("I88").Value = 2
$ws.Range("I89").Value = 0
$ws.Range("J88").Value = 2
$ws.Range("J89").Value = 0
$ws.Range("K88").Value = 'A'
$ws.Range("K89").Value = 'D'
$ws.Range("L88").Value = 3.6
$ws.Range("L89").Value = 1.615
$ws.Range("M88").Value = 3.4
$ws.Range("M89").Value = 4
$ws.Range("N88").Value = 1.909
$ws.Range("N89").Value = 4.5
$ws.Range("O88").Value = 2.4
$ws.Range("O89").Value = 1.85
$ws.Range("P88").Value = 3.6
$ws.Range("P89").Value = 3.8
$ws.Range("Q88").Value = 2.45
$ws.Range("Q89").Value = 3.5
$ws.Range("R88").Value = 0
$ws.Range("R89").Value = -0.5
$ws.Range("S88").Value = 1.875
$ws.Range("S89").Value = 1.875
$ws.Range("T88").Value = 1.925
$ws.Range("T89").Value = 1.925
$ws.Range("U88").Value = 2.75
$ws.Range("U89").Value = 2.5
$ws.Range("V88").Value = 1.975
$ws.Range("V89").Value = 1.75
$ws.Range("W88").Value = 1.825
$ws.Range("W89").Value = 1.95
$ws.Range("X88").Value = -1
$ws.Range("X89").Value = -1
$ws.Range("Y88").Value = -1
$ws.Range("Y89").Value = 2.8
$ws.Range("Z88").Value = 1.45
$ws.Range("Z89").Value = -1
$ws.Range("AA88").Value = -1
$ws.Range("AA89").Value = -1
$ws.Range("AB88").Value = 0.925
$ws.Range("AB89").Value = 0.925
$ws.Range("AC88").Value = 0.9750000000000001
$ws.Range("AC89").Value = -1
$ws.Range("AD88").Value = -1
$ws.Range("AD89").Value = 0.95

# --- Swap row 104 and row 106 (columns B, E:AD); A/C/D unchanged ---
$ws.Range("B104").Value = 6533597
$ws.Range("B106").Value = 6537869
$ws.Range("E104").Value = 'FC Kuressaare'
$ws.Range("E106").Value = 'JK Tallinna Kalev'
$ws.Range("F104").Value = 'Parnu JK Vaprus'
$ws.Range("F106").Value = 'JK Trans Narva'
$ws.Range("G104").Value = 1
$ws.Range("G106").Value = 5
$ws.Range("H104").Value = 0
$ws.Range("H106").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("I106").Value = 2
$ws.Range("J104").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K104").Value = 'H'
$ws.Range("K106").Value = 'H'
$ws.Range("L104").Value = 2.5
$ws.Range("L106").Value = 1.6
$ws.Range("M104").Value = 3.4
$ws.Range("M106").Value = 4
$ws.Range("N104").Value = 2.5
$ws.Range("N106").Value = 4.5
$ws.Range("O104").Value = 2.15
$ws.Range("O106").Value = 1.65
$ws.Range("P104").Value = 3.6
$ws.Range("P106").Value = 4
$ws.Range("Q104").Value = 2.875
$ws.Range("Q106").Value = 4.333
$ws.Range("R104").Value = -0.25
$ws.Range("R106").Value = -0.75
$ws.Range("S104").Value = 1.95
$ws.Range("S106").Value = 1.8
$ws.Range("T104").Value = 1.85
$ws.Range("T106").Value = 2
$ws.Range("U104").Value = 2.75
$ws.Range("U106").Value = 2.75
$ws.Range("V104").Value = 1.95
$ws.Range("V106").Value = 1.9
$ws.Range("W104").Value = 1.85
$ws.Range("W106").Value = 1.9
$ws.Range("X104").Value = 1.15
$ws.Range("X106").Value = 0.6499999999999999
$ws.Range("Y104").Value = -1
$ws.Range("Y106").Value = -1
$ws.Range("Z104").Value = -1
$ws.Range("Z106").Value = -1
$ws.Range("AA104").Value = 0.95
$ws.Range("AA106").Value = 0.8
$ws.Range("AB104").Value = -1
$ws.Range("AB106").Value = -1
$ws.Range("AC104").Value = -1
$ws.Range("AC106").Value = 0.8999999999999999
$ws.Range("AD104").Value = 0.8500000000000001
$ws.Range("AD106").Value = -1

# --- Swap row 120 and row 121 (columns B, E:AD); A/C/D unchanged ---
$ws.Range("B120").Value = 7721007
$ws.Range("B121").Value = 7721087
$ws.Range("E120").Value = 'JK Trans Narva'
$ws.Range("E121").Value = 'Paide Linnameeskond'
$ws.Range("F120").Value = 'JK Tammeka Tartu'
$ws.Range("F121").Value = 'FC Flora Tallinn'
$ws.Range("G120").Value = 0
$ws.Range("G121").Value = 2
$ws.Range("H120").Value = 5
$ws.Range("H121").Value = 1
$ws.Range("I120").Value = 0
$ws.Range("I121").Value = 1
$ws.Range("J120").Value = 2
$ws.Range("J121").Value = 1
$ws.Range("K120").Value = 'A'
$ws.Range("K121").Value = 'H'
$ws.Range("L120").Value = 2.25
$ws.Range("L121").Value = 2.2
$ws.Range("M120").Value = 3.3
$ws.Range("M121").Value = 3.3
$ws.Range("N120").Value = 2.75
$ws.Range("N121").Value = 2.8
$ws.Range("O120").Value = 2.1
$ws.Range("O121").Value = 1.85
$ws.Range("P120").Value = 3.25
$ws.Range("P121").Value = 3.6
$ws.Range("Q120").Value = 3
$ws.Range("Q121").Value = 3.4
$ws.Range("R120").Value = -0.25
$ws.Range("R121").Value = -0.5
$ws.Range("S120").Value = 1.875
$ws.Range("S121").Value = 1.9
$ws.Range("T120").Value = 1.925
$ws.Range("T121").Value = 1.9
$ws.Range("U120").Value = 2.5
$ws.Range("U121").Value = 2.5
$ws.Range("V120").Value = 1.825
$ws.Range("V121").Value = 1.95
$ws.Range("W120").Value = 1.975
$ws.Range("W121").Value = 1.85
$ws.Range("X120").Value = -1
$ws.Range("X121").Value = 0.8500000000000001
$ws.Range("Y120").Value = -1
$ws.Range("Y121").Value = -1
$ws.Range("Z120").Value = 2
$ws.Range("Z121").Value = -1
$ws.Range("AA120").Value = -1
$ws.Range("AA121").Value = 0.8999999999999999
$ws.Range("AB120").Value = 0.925
$ws.Range("AB121").Value = -1
$ws.Range("AC120").Value = 0.825
$ws.Range("AC121").Value = 0.95
$ws.Range("AD120").Value = -1
$ws.Range("AD121").Value = -1

# --- Rebuild the final match row: old row 176 data shifts down to row 178,
#     two new match rows are inserted as 176 and 177 ---

# New row 176 (new match: id 174)
$ws.Range("A176").Value = 174
$ws.Range("B2").Copy()
$ws.Range("A176").PasteSpecial(-4122)
$ws.Range("B176").Value = 7719660
$ws.Range("C176").Value = "Estonia Meistriliiga"
$ws.Range("D176").Value = 45441.54166666666
$ws.Range("D2").Copy()
$ws.Range("D176").PasteSpecial(-4122)
$ws.Range("E176").Value = "JK Tallinna Kalev"
$ws.Range("F176").Value = "FC Flora Tallinn"
$ws.Range("G176").Value = 2
$ws.Range("H176").Value = 3
$ws.Range("K176").Value = "A"
$ws.Range("L176").Value = 5
$ws.Range("M176").Value = 4.2
$ws.Range("N176").Value = 1.5
$ws.Range("O176").Value = 3.9
$ws.Range("P176").Value = 3.9
$ws.Range("Q176").Value = 1.727
$ws.Range("R176").Value = 0.75
$ws.Range("S176").Value = 1.85
$ws.Range("T176").Value = 1.95
$ws.Range("U176").Value = 3
$ws.Range("V176").Value = 1.85
$ws.Range("W176").Value = 1.95
$ws.Range("X176").Value = -1
$ws.Range("Y176").Value = -1
$ws.Range("Z176").Value = 0.7270000000000001
$ws.Range("AA176").Value = -0.5
$ws.Range("AB176").Value = 0.475
$ws.Range("AC176").Value = 0.8500000000000001
$ws.Range("AD176").Value = -1

# New row 177 (new match: id 175)
$ws.Range("A177").Value = 175
$ws.Range("B2").Copy()
$ws.Range("A177").PasteSpecial(-4122)
$ws.Range("B177").Value = 7721036
$ws.Range("C177").Value = "Estonia Meistriliiga"
$ws.Range("D177").Value = 45441.54166666666
$ws.Range("D2").Copy()
$ws.Range("D177").PasteSpecial(-4122)
$ws.Range("E177").Value = "Paide Linnameeskond"
$ws.Range("F177").Value = "FC Levadia Tallinn"
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 1
$ws.Range("K177").Value = "A"
$ws.Range("L177").Value = 4.5
$ws.Range("M177").Value = 3.8
$ws.Range("N177").Value = 1.615
$ws.Range("O177").Value = 4.1
$ws.Range("P177").Value = 3.8
$ws.Range("Q177").Value = 1.65
$ws.Range("R177").Value = 0.75
$ws.Range("S177").Value = 1.925
$ws.Range("T177").Value = 1.875
$ws.Range("U177").Value = 2.75
$ws.Range("V177").Value = 1.8
$ws.Range("W177").Value = 2
$ws.Range("X177").Value = -1
$ws.Range("Y177").Value = -1
$ws.Range("Z177").Value = 0.6499999999999999
$ws.Range("AA177").Value = -0.5
$ws.Range("AB177").Value = 0.4375
$ws.Range("AC177").Value = -1
$ws.Range("AD177").Value = 1

# Row 178 = the match formerly stored at row 176 (id 176), with O adjusted 3.25 -> 3.2
$ws.Range("A178").Value = 176
$ws.Range("B2").Copy()
$ws.Range("A178").PasteSpecial(-4122)
$ws.Range("B178").Value = "7721037"
$ws.Range("C178").Value = "Estonia Meistriliiga"
$ws.Range("D178").Value = 45443.58333333334
$ws.Range("D2").Copy()
$ws.Range("D178").PasteSpecial(-4122)
$ws.Range("E178").Value = "JK Nomme United"
$ws.Range("F178").Value = "Parnu JK Vaprus"
$ws.Range("L178").Value = 3.5
$ws.Range("M178").Value = 3.4
$ws.Range("N178").Value = 1.909
$ws.Range("O178").Value = 3.2
$ws.Range("P178").Value = 3.4
$ws.Range("Q178").Value = 2
$ws.Range("R178").Value = 0.25
$ws.Range("S178").Value = 2
$ws.Range("T178").Value = 1.8
$ws.Range("U178").Value = 2.5
$ws.Range("V178").Value = 1.85
$ws.Range("W178").Value = 1.95
$ws.Range("X178").Value = 0
$ws.Range("Y178").Value = 0
$ws.Range("Z178").Value = 0
